# Positions.xlsx update:
#  1) Existing AB column totals (rows 2-524) are rescaled by 4/3.
#  2) A new data row (525) is appended with the latest snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rescale AB2:AB524 by 4/3 ---------------------------------------
for ($r = 2; $r -le 524; $r++) {
    $cell = $ws.Cells.Item($r, 28)   # column AB = 28
    $cell.Value = $cell.Value2 * 4 / 3
}

# --- 2) Append new row 525 ----------------------------------------------
$row = 525
$ws.Cells.Item($row, 1).Value = 45955              # A - date serial

# Match column-A formatting (bold, bordered, centered, custom date numFmt)
# used by every other row in the sheet.
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)        # xlPasteFormats

$ws.Cells.Item($row, 2).Value  = 1359.5668706115      # B  BTCUSDT
$ws.Cells.Item($row, 3).Value  = 0.143721357          # C  ETHUSDT
$ws.Cells.Item($row, 4).Value  = 0                    # D  FDUSDUSDT
$ws.Cells.Item($row, 5).Value  = 0                    # E  FETUSDT
$ws.Cells.Item($row, 6).Value  = 0                    # F  INJUSDT
$ws.Cells.Item($row, 7).Value  = 0                    # G  MKRUSDT
$ws.Cells.Item($row, 8).Value  = 0                    # H  RNDRUSDT
$ws.Cells.Item($row, 9).Value  = 0.08589050399999999  # I  SOLUSDT
$ws.Cells.Item($row, 10).Value = 0                    # J  TNSRUSDT
$ws.Cells.Item($row, 11).Value = 575.746628556233     # K  TRXUSDT
$ws.Cells.Item($row, 12).Value = 0                    # L  UMAUSDT
$ws.Cells.Item($row, 13).Value = 0                    # M  USDTUSDT
$ws.Cells.Item($row, 14).Value = 0.0015487904         # N  BBUSDT
$ws.Cells.Item($row, 15).Value = 0.0293870696         # O  BNBUSDT
$ws.Cells.Item($row, 16).Value = 0                    # P  NEARUSDT
$ws.Cells.Item($row, 17).Value = 0.000001716          # Q  PEPEUSDT
$ws.Cells.Item($row, 18).Value = 0                    # R  XRPUSDT
$ws.Cells.Item($row, 19).Value = 0                    # S  ARUSDT
$ws.Cells.Item($row, 20).Value = 0                    # T  ENSUSDT
$ws.Cells.Item($row, 21).Value = 191.6376691618858    # U  JASMYUSDT
$ws.Cells.Item($row, 22).Value = 0                    # V  LDUMAUSDT
$ws.Cells.Item($row, 23).Value = 0                    # W  LPTUSDT
$ws.Cells.Item($row, 24).Value = 0                    # X  NOTUSDT
$ws.Cells.Item($row, 25).Value = 0                    # Y  EURUSDT
$ws.Cells.Item($row, 26).Value = 54.079996785894      # Z  OMUSDT
$ws.Cells.Item($row, 27).Value = 0.322727731542       # AA USDCUSDT
$ws.Cells.Item($row, 28).Value = 8743.629713818566    # AB Total

Write-Output "done"
